# Update Fn1-Itga5 LR-pair sheet with newly recomputed TPM-derived expression
# values. Ligand-side stats (G,H,I,J) depend only on the Sending cluster;
# receptor-side stats (M,N,O,P) depend only on the Target cluster; the edge
# stats (Q,R,S,T) are combinations of the two and are set per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-expression stats, keyed by Sending cluster.
# [Ligand avg expr, Ligand total expr, Ligand specificity (avg), Ligand specificity (total)]
$ligandBySender = @{
    "ECs"               = @(29.20950566666667, 87.628517,   0.01829497698069002, 0.01840828041918582)
    "FAPs"              = @($null,             $null,       0.913374480506715,   0.9190311407684336)
    "Inflammatory-Mac"  = @(57.98602933333333, 173.958088,  0.03631876156896331, 0.03654368891224535)
    "MuSCs"             = @(29.481085,         58.96217,    0.01846507700595112, 0.01238628926567028)
    "Resolving-Mac"     = @(21.628479,         64.885437,   0.01354670393768061, 0.01363060063446486)
}

# New receptor-expression stats, keyed by Target cluster.
# [Receptor avg expr, Receptor total expr, Receptor specificity (avg), Receptor specificity (total)]
$receptorByTarget = @{
    "ECs"               = @(28.85518433333334,  86.56555300000001, 0.1999651185353207, 0.2044513327926365)
    "FAPs"              = @($null,               $null,            0.3546352265743414, 0.3625914622481308)
    "Inflammatory-Mac"  = @(29.393479,           88.180437,         0.2036954761578358, 0.2082653809291453)
    "MuSCs"             = @(9.499066500000001,   18.998133,         0.0658280999596015, 0.04486996822421697)
    "Resolving-Mac"     = @(25.37910966666666,   76.13732899999999, 0.1758760787729007, 0.1798218558058706)
}

# New edge-weight stats, keyed by row number.
# [Edge avg weight, Edge total weight, Edge avg-derived specificity, Edge total-derived specificity]
$edgeByRow = @{
    2  = @(842.8456702972114,  7585.611032674902,  0.003658357240544643, 0.003763597466123134)
    3  = @(1494.774525889413,  13452.97073300472,  0.006488043306719366, 0.006674685314666222)
    4  = @(858.5689914135477,  7727.120922721929,  0.003726604047378299, 0.003833807533752261)
    5  = @(277.4630367597936,  1664.778220558761,  0.001204323573443471, 0.0008259789574713432)
    6  = @(741.3112476234547,  6671.801228611092,  0.003217648812604242, 0.003310211147172864)
    7  = @(42078.96665120669,  378710.6998608602,  0.1826430362616622,   0.1878971416080434)
    8  = @(74626.43475856508,  671637.9128270858,  0.3239147658417203,   0.3332328451827941)
    9  = @(42863.95152829246,  385775.5637546322,  0.1860502497172313,   0.1914023706178848)
    10 = @(13852.30805853321,  83113.8483511993,   0.06012570660334512,  0.04123689808334548)
    11 = @(37009.87306004909,  333088.8575404418,  0.1606407220827562,   0.1652618852763661)
    12 = @(1673.197565171407,  15058.77808654266,  0.007262485462193798, 0.007471405903268054)
    13 = @(2967.39152295398,   26706.52370658582,  0.01287991223790879,  0.01325042959863184)
    14 = @(1704.411135502717,  15339.70021952445,  0.007397967431252891, 0.007610785291864959)
    15 = @(550.8131487082841,  3304.878892249704,  0.00239079506697065,  0.001639714160288119)
    16 = @(1471.633797585217,  13244.70417826695,  0.006387601370637188, 0.006571353958192376)
    17 = @(850.6821420216685,  5104.09285213001,   0.003692371312258841, 0.002532393348721415)
    18 = @(1508.672394407193,  9052.03436644316,   0.006548366767718139, 0.004491162736667713)
    19 = @(866.551652844715,   5199.30991706829,   0.00376125265301832,  0.002579635252213403)
    20 = @(280.0427869071526,  1120.17114762861,   0.00121552093490949,  0.0005557724057665851)
    21 = @(748.2036893073216,  4489.22213584393,   0.003247565338046337, 0.002227325522301164)
    22 = @(624.0937483946291,  5616.843735551661,  0.002708868258661198, 0.002786794466480498)
    23 = @(1106.821177046764,  9961.390593420876,  0.004804138420274886, 0.004942339415370914)
    24 = @(635.736243288441,   5721.626189595969,  0.002759402308955082, 0.002838782233429874)
    25 = @(205.4503603148535,  1232.702161889121,  0.0008917537809327664,0.0006116046173454301)
    26 = @(548.9115404641969,  4940.203864177773,  0.002382541168856678, 0.002451079901838149)
}

$firstRow = 2
$lastRow = 26

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $sender = $ws.Cells.Item($row, 1).Value2   # column A: Sending cluster
    $target = $ws.Cells.Item($row, 4).Value2   # column D: Target cluster

    $lig = $ligandBySender[$sender]
    if ($lig[0] -ne $null) { $ws.Cells.Item($row, 7).Value  = $lig[0] }   # G: Ligand average expression value
    if ($lig[1] -ne $null) { $ws.Cells.Item($row, 8).Value  = $lig[1] }   # H: Ligand total expression value
    $ws.Cells.Item($row, 9).Value  = $lig[2]                              # I: Ligand derived specificity (avg)
    $ws.Cells.Item($row, 10).Value = $lig[3]                              # J: Ligand derived specificity (total)

    $rec = $receptorByTarget[$target]
    if ($rec[0] -ne $null) { $ws.Cells.Item($row, 13).Value = $rec[0] }  # M: Receptor average expression value
    if ($rec[1] -ne $null) { $ws.Cells.Item($row, 14).Value = $rec[1] }  # N: Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $rec[2]                              # O: Receptor derived specificity (avg)
    $ws.Cells.Item($row, 16).Value = $rec[3]                              # P: Receptor derived specificity (total)

    $edge = $edgeByRow[$row]
    $ws.Cells.Item($row, 17).Value = $edge[0]   # Q: Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $edge[1]   # R: Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $edge[2]   # S: Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $edge[3]   # T: Edge total expression derived specificity
}
